# CACEQUI.xlsx automatic update
#
# Changes applied (per the commit's OOXML diff):
#   1. Rename sheet "Paineis DARQ"            -> "PAINEIS DARQ"
#   2. Rename sheet "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   3. Remove the sheet "Desarquivamentos Pendentes" entirely
#      (its two unique strings "PEDIDOS PENDENTES"/"%" and its five
#      unique cell-formats were only referenced by that sheet, so Excel
#      drops them from sharedStrings.xml/styles.xml automatically once
#      the sheet disappears)
#   4. The "DGC" sheet's tab itself is untouched content-wise - it just
#      shifts up one slot in the workbook now that the preceding sheet
#      is gone.

$wb = $excel.ActiveWorkbook

# 1) Rename "Paineis DARQ" -> "PAINEIS DARQ"
$painelSheet = $wb.Worksheets.Item("Paineis DARQ")
$painelSheet.Name = "PAINEIS DARQ"

# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$recolhimentoSheet = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolhimentoSheet.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3) Delete the "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$pendentesSheet = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$pendentesSheet.Delete()
$excel.DisplayAlerts = $true
